$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.888795852661133
$ws.Range("B1").Value = 3.51558780670166
$ws.Range("C1").Value = 3.087420701980591
$ws.Range("D1").Value = 3.353963851928711
$ws.Range("E1").Value = 1.875334501266479
